# Update cryptocurrency price/volume data for Wed Apr  5 02:47:36 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep the Price column (D) as text so values like "1.000" or "0.06932"
# are not reinterpreted by Excel as numbers and lose their original formatting.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.588.83"
$ws.Range("E2").Value = "  +2.30%  "
$ws.Range("D3").Value = "1.911.19"
$ws.Range("E3").Value = "  +5.43%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "313.44"
$ws.Range("E5").Value = "  +1.09%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").Value = "0.5036"
$ws.Range("E7").Value = "  +2.00%  "
$ws.Range("D8").Value = "0.3959"
$ws.Range("E8").Value = "  +2.19%  "
$ws.Range("D9").Value = "0.09642"
$ws.Range("E9").Value = "  -1.87%  "
$ws.Range("D10").Value = "1.159"
$ws.Range("E10").Value = "  +5.20%  "
$ws.Range("D11").Value = "41.64"
$ws.Range("E11").Value = "  +1.89%  "
$ws.Range("D12").Value = "6.564"
$ws.Range("E12").Value = "  +1.95%  "
$ws.Range("D13").Value = "21.13"
$ws.Range("E13").Value = "  +2.96%  "
$ws.Range("D14").Value = "1.912.97"
$ws.Range("E14").Value = "  +5.58%  "
$ws.Range("D15").Value = "7.557"
$ws.Range("E15").Value = "  +3.50%  "
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").Value = "  +0.20%  "
$ws.Range("D17").Value = "0.00001135"
$ws.Range("E17").Value = "  -0.07%  "
$ws.Range("D18").Value = "93.89"
$ws.Range("E18").Value = "  +1.28%  "
$ws.Range("E19").Value = "  +0.21%  "
$ws.Range("D20").Value = "18.05"
$ws.Range("E20").Value = "  +5.58%  "
$ws.Range("D21").Value = "0.9998"
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("D22").Value = "6.268"
$ws.Range("E22").Value = "  +5.54%  "
$ws.Range("D23").Value = "28.664.20"
$ws.Range("E23").Value = "  +2.43%  "
$ws.Range("D24").Value = "11.43"
$ws.Range("E24").Value = "  +2.55%  "
$ws.Range("D25").Value = "2.283"
$ws.Range("E25").Value = "  +1.67%  "
$ws.Range("D26").Value = "2.765"
$ws.Range("E26").Value = "  +15.60%  "
$ws.Range("D27").Value = "2.140.34"
$ws.Range("E27").Value = "  +5.91%  "
$ws.Range("D28").Value = "21.39"
$ws.Range("E28").Value = "  +3.95%  "
$ws.Range("D29").Value = "159.20"
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("D30").Value = "128.74"
$ws.Range("E30").Value = "  +1.08%  "
$ws.Range("D31").Value = "1.115"
$ws.Range("E31").Value = "  +7.12%  "
$ws.Range("D32").Value = "0.1076"
$ws.Range("E32").Value = "  +1.61%  "
$ws.Range("D33").Value = "5.719"
$ws.Range("E33").Value = "  +2.49%  "
$ws.Range("D34").Value = "3.631"
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("D35").Value = "9.801"
$ws.Range("E35").Value = "  +8.69%  "
$ws.Range("D36").Value = "0.06787"
$ws.Range("E36").Value = "  +0.41%  "
$ws.Range("D37").Value = "0.02439"
$ws.Range("E37").Value = "  +4.77%  "
$ws.Range("D38").Value = "0.2220"
$ws.Range("E38").Value = "  +4.01%  "
$ws.Range("D39").Value = "5.105"
$ws.Range("E39").Value = "  +3.31%  "
$ws.Range("D40").Value = "11.62"
$ws.Range("E40").Value = "  +3.22%  "
$ws.Range("D41").Value = "0.6413"
$ws.Range("E41").Value = "  +3.44%  "
$ws.Range("D42").Value = "1.197"
$ws.Range("E42").Value = "  +4.52%  "
$ws.Range("D43").Value = "0.9990"
$ws.Range("E43").Value = "  -0.07%  "
$ws.Range("D44").Value = "13.73"
$ws.Range("E44").Value = "  +4.42%  "
$ws.Range("D45").Value = "0.6091"
$ws.Range("E45").Value = "  +3.74%  "
$ws.Range("E47").Value = "  -0.96%  "
$ws.Range("D48").Value = "2.040"
$ws.Range("E48").Value = "  +5.53%  "
$ws.Range("D49").Value = "124.92"
$ws.Range("E49").Value = "  +1.95%  "
$ws.Range("D50").Value = "1.209"
$ws.Range("E50").Value = "  +2.91%  "

# Row 51: Aave dropped out of the top list and was replaced by Cronos
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.06932"
$ws.Range("E51").Value = "  +2.17%  "
